# Update the violations dictionary sheet:
#  - translate the Russian "7 дней" values in column J (time_to_correct)
#    to the English "7 days" for the listed rows
#  - fill in the previously empty J63 / J65 cells with "5 minutes"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J is index 10 ("time_to_correct")
$colJ = 10

$rowsToEnglish = @(10, 12, 18, 24, 31, 44, 46, 52, 54, 56, 62, 70, 76, 85, 87, 91, 95, 99, 103, 105, 107)
foreach ($r in $rowsToEnglish) {
    $ws.Cells.Item($r, $colJ).Value = "7 days"
}

$rowsFiveMinutes = @(63, 65)
foreach ($r in $rowsFiveMinutes) {
    $ws.Cells.Item($r, $colJ).Value = "5 minutes"
}
